# Updated Class Responsibilities and Design Rationale for Beating up the
# Zombies: expand the terse "Execute an action" / "The action is the first
# non-null action ..." bullets under the Zombie heading into their final
# wording, and add two new bullets ("Keep track of the number of arms and
# legs it has" and "Create and drop a WeaponItem when a limb is lost").

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: build a minimal Flat-OPC "pkg:package" wrapper around a body
# fragment so Range.InsertXML can splice in exact <w:r>/<w:p> runs
# without Word's usual "merge adjacent identically-formatted runs"
# behaviour collapsing the pieces back together.
# ---------------------------------------------------------------------
function New-FlatOpcBody([string]$bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------
# 1) "Execute an action"
#    -> "Execute the first non-null action returned by ScavengeBehaviour,
#        AttackBehaviour, HuntBeviour or WanderBehaviour, in that order"
#    (split across six runs, matching how the final sentence's tail was
#    carried over verbatim from the paragraph below)
# ---------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("Execute an action", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find paragraph 'Execute an action'"
}
$firstPara = $find1.Paragraphs(1).Range
$firstPara.End = $firstPara.End - 1
$firstPara.Delete()

$run1Body = '<w:p><w:r><w:t xml:space="preserve">Execute </w:t></w:r><w:r><w:t xml:space="preserve">the first non-null action </w:t></w:r><w:r><w:t>returned by ScavengeBehaviour, AttackBehaviour, HuntBeviour</w:t></w:r><w:r><w:t xml:space="preserve"> or</w:t></w:r><w:r><w:t xml:space="preserve"> WanderBehaviour</w:t></w:r><w:r><w:t>, in that order</w:t></w:r></w:p>'
$ins1 = $d.Range($firstPara.Start, $firstPara.Start)
$ins1.InsertXML((New-FlatOpcBody $run1Body))

# ---------------------------------------------------------------------
# 2) "The action is the first non-null action returned by
#    ScavengeBehaviour, AttackBehaviour, HuntBeviour or WanderBehaviour,
#    in that order" (plus the trailing _GoBack bookmark) becomes three
#    things:
#      - "Keep track of the number of arms and legs it has"
#      - "Create and drop a WeaponItem when a limb is lost" (new bullet,
#        now carrying the _GoBack bookmark at its end)
# ---------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("The action is the first non-null action", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find paragraph 'The action is the first non-null action ...'"
}
$secondPara = $find2.Paragraphs(1).Range
$secondPara.End = $secondPara.End - 1
$secondPara.Delete()

$run2Body = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Keep track of the number of arms and legs it has</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Create </w:t></w:r><w:r><w:t xml:space="preserve">and drop </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/></w:rPr><w:t>WeaponItem</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> when a limb is lost</w:t></w:r></w:p>
'@
$ins2 = $d.Range($secondPara.Start, $secondPara.Start)
$ins2.InsertXML((New-FlatOpcBody $run2Body))

Write-Host "Done."
